# Change self data item template handling
#
# Inserts a new numeric "self" value column (D) in front of the existing
# label column, pushing the old label column from D to E for the
# Test1/Test2/Test3 total rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing D column (self-item template labels) one slot to the
# right, making room for the new numeric self-item value column.
$ws.Columns("D").Insert()

# Populate the new D column with the self data item's numeric value for
# each of the rows that previously only held the template label.
$ws.Range("D5").Value = 55.76
$ws.Range("D9").Value = 5500.8

# Match the resulting cursor/selection position left behind by the edit.
$ws.Range("F11").Select()
